$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# --- D/E columns: file name references (TC01/DNBSEQ-G400 -> TC09/Illumina NextSeq) ---
$dVal = "TC09_CDS_Filter_InstrumentModel-Illumina NextSeq_Neo4jData.xlsx"
$eVal = "TC09_CDS_Filter_InstrumentModel-Illumina NextSeq_WebData.xlsx"

$ws.Range("D2").Value = $dVal
$ws.Range("D3").Value = $dVal
$ws.Range("D4").Value = $dVal

$ws.Range("E2").Value = $eVal
$ws.Range("E3").Value = $eVal
$ws.Range("E4").Value = $eVal

# --- C column: clear old "Files" query text first so its shared-string slot is
#     freed before the B-column strings are rewritten (matches the shared
#     string table ordering produced by Excel when re-entering this content) ---
$ws.Range("C2").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("C4").ClearContents()

# --- B2: Participants query ---
$b2Val = "Match (f)<--(g:genomic_info)`nWHERE g.instrument_model in ['Illumina NextSeq']`nMATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)`nWITH p, s, collect(distinct samp.sample_id) as samp`nRETURN `ncoalesce(p.participant_id,'') as ``Participant ID``,`ncoalesce(s.study_name, '') as ``Study Name``,`ncoalesce(s.phs_accession,'') as ``Accession``,`ncoalesce(p.gender,'') as ``Gender``,`ncoalesce(apoc.text.join(samp, ','), '') as ``Samples```nORDER BY ``Participant ID``LIMIT 100"
$ws.Range("B2").Value = $b2Val

# --- B3: Samples query ---
$b3Val = "Match (f)<--(g:genomic_info)`nWHERE g.instrument_model in ['Illumina NextSeq']`nMATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)`nWITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor`nRETURN  `n coalesce(samp.sample_id, '') as ``Sample ID``,`n coalesce(p.participant_id,'') as ``Participant ID``,`n coalesce(s.study_name, '') as ``Study Name``,`n coalesce(s.phs_accession,'') as ``Accession``,`ncoalesce(samp.sample_tumor_status,'') as ``Tumor``,`ncoalesce(samp.sample_type,'') as ``Analyte Type```nORDER By samp.sample_id LIMIT 100"
$ws.Range("B3").Value = $b3Val

# --- B4: Files (names) query -- note WHERE line is joined directly to MATCH line (no newline) ---
$b4Val = "Match (f)<--(g:genomic_info)`nWHERE g.instrument_model in ['Illumina NextSeq']MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)`nWITH p,s,f,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor`nRETURN `n    coalesce(f.file_name, '') as ``File Name``,`n    coalesce(s.study_name, '') as ``Study Name``,`n    coalesce(s.phs_accession,'') as ``Accession``,`n    coalesce(p.participant_id,'') as ``Participant ID``,`n    coalesce(samp.sample_id, '') as ``Sample ID``,`n    coalesce(f.file_type, '') as ``File Type```nORDER By f.file_name LIMIT 100"
$ws.Range("B4").Value = $b4Val

# --- C column: "Files" Cypher query (StatQuery), instrument_model filter updated;
#     re-entered after the B-column edits so it lands as a freshly appended
#     shared string (matches the table ordering captured in the diff) ---
$cVal = "MATCH (f:file)`nMatch (f)<--(g:genomic_info)`nWHERE g.instrument_model in ['Illumina NextSeq']`nMATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)`nWITH p,f, s, collect(distinct samp.sample_id) as samp`nRETURN`ncount(distinct s) AS Studies,`ncount(distinct p) AS Participants,`ncount(distinct samp) AS Samples,`ncount(distinct f) AS Files"

$ws.Range("C2").Value = $cVal
$ws.Range("C3").Value = $cVal
$ws.Range("C4").Value = $cVal

# --- column D widens (bestFit) to accommodate the longer "Illumina NextSeq"
#     file-name text now stored in D2:D4 ---
$ws.Columns.Item(4).ColumnWidth = 90.8

# --- update selection to D3 (matches saved sheetView selection in diff) ---
$ws.Range("D3").Select()
